# Insert two new weekly price rows (date 2021-12-29, serial 44559) at the
# top of the "Pepino ensalada" data block for Agrícola del Norte S.A. de
# Arica, pushing the existing rows 171-227 down to 173-229.

$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

# Insert two blank rows above the current row 171 (data rows 171-227 shift
# down to 173-229; formats/styles are inherited from the row below, as
# Excel does natively).
$ws.Rows("171:172").Insert()

# New row 171: "Primera" quality, $/caja 70 unidades
$ws.Range("A171").Value = 1
$ws.Range("B171").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C171").Value = "Arica y Parinacota"
$ws.Range("D171").Value = 44559
$ws.Range("E171").Value = 15
$ws.Range("F171").Value = 100112043
$ws.Range("G171").Value = "Pepino ensalada"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 120
$ws.Range("K171").Value = 3500
$ws.Range("L171").Value = 4000
$ws.Range("M171").Value = 3750
$ws.Range("N171").Value = "$/caja 70 unidades"
$ws.Range("O171").Value = "Región de Arica y Parinacota"
$ws.Range("P171").Value = 54
$ws.Range("Q171").Value = 70
$ws.Range("R171").Value = "Hortaliza"

# New row 172: "Segunda" quality, $/caja 100 unidades
$ws.Range("A172").Value = 1
$ws.Range("B172").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C172").Value = "Arica y Parinacota"
$ws.Range("D172").Value = 44559
$ws.Range("E172").Value = 15
$ws.Range("F172").Value = 100112043
$ws.Range("G172").Value = "Pepino ensalada"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Segunda"
$ws.Range("J172").Value = 120
$ws.Range("K172").Value = 2500
$ws.Range("L172").Value = 3000
$ws.Range("M172").Value = 2750
$ws.Range("N172").Value = "$/caja 100 unidades"
$ws.Range("O172").Value = "Región de Arica y Parinacota"
$ws.Range("P172").Value = 28
$ws.Range("Q172").Value = 100
$ws.Range("R172").Value = "Hortaliza"
